# rev1.6: add L6 (optional) inductor row between L3 (row 23) and L2 (row 24),
# and update the sheet view (zoom + selection) and workbook tab-ratio to
# match the new revision's saved state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# Insert a new row above the current row 24 (the "L2 / 20n" row), pushing
# L2 and everything below it down by one row.
$ws.Rows("24:24").Insert()

# Populate the new row 24 with the L6 (18n, optional) inductor data -
# same footprint / distributor / MFR as the neighbouring inductors.
$ws.Range("A24").Value = "L6"
$ws.Range("B24").Value = "18n"
$ws.Range("C24").Value = "Inductors_SMD:L_0402"
$ws.Range("D24").Value = 1
$ws.Range("E24").Value = "TE Connectivity / Sigma Inductors"
$ws.Range("F24").Value = "36501E18NJTDG"
$ws.Range("G24").Value = "Mouser"
$ws.Range("H24").Value = "279-36501E18NJTDG"

# Update the saved view: selection moved to A24, zoom reset to 100%.
$ws.Range("A24").Select()
$excel.ActiveWindow.Zoom = 100

# Workbook window tab-ratio (sheet-tab/scrollbar split) changed in this
# revision.
$excel.ActiveWindow.TabRatio = 0.5

# Print setup: paper size changed to A4 (9) in this revision.
$ws.PageSetup.PaperSize = 9
